# Check for blank exchange rates for import distributions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the exchange-rate related cells (From Currency, To Currency, Exchange Rate,
# Commitment Percentage) on row 4, but keep the "As Of" date cell present
# (blank, retaining its date number-format).
$ws.Range("L4:O4").ClearContents()

# Move the active selection to O4 (matches the post-edit selection in the file)
$ws.Range("O4").Select()
